$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1955.3
$ws.Range("I28").Value = 1839.2222
$ws.Range("J28").Value = 3000
$ws.Range("K28").Value = 1839.2222
$ws.Range("L28").Value = 3000
$ws.Range("M28").Value = -1354.2222
$ws.Range("N28").Value = -3970

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 49.5
$ws.Range("I39").Value = 45.857143
$ws.Range("J39").Value = 75
$ws.Range("K39").Value = 137.571429
$ws.Range("L39").Value = 225
$ws.Range("M39").Value = 158.428571
$ws.Range("N39").Value = -817

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").ClearContents()
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = 0

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").ClearContents()
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = 0

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 4168
$ws.Range("I101").Value = 4168
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 12504
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = -10882

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 28065.5
$ws.Range("I113").Value = 21678.6
$ws.Range("J113").Value = 60000
$ws.Range("K113").Value = 21678.6
$ws.Range("L113").Value = 60000
$ws.Range("M113").Value = -18424.6
$ws.Range("N113").Value = -66508

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 6799.4
$ws.Range("I138").Value = 2000
$ws.Range("J138").Value = 9999
$ws.Range("K138").Value = 6000
$ws.Range("L138").Value = 29997
$ws.Range("M138").Value = -860
$ws.Range("N138").Value = -40277

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1948.5
$ws.Range("I45").Value = 1948.5
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1948.5
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -1571.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2186.25
$ws.Range("I88").Value = 1997.5
$ws.Range("J88").Value = 2249.1667
$ws.Range("K88").Value = 1997.5
$ws.Range("L88").Value = 2249.1667
$ws.Range("M88").Value = -1591.5
$ws.Range("N88").Value = -3061.1667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 2186.25
$ws.Range("I91").Value = 1997.5
$ws.Range("J91").Value = 2249.1667
$ws.Range("K91").Value = 1997.5
$ws.Range("L91").Value = 2249.1667
$ws.Range("M91").Value = -593.5
$ws.Range("N91").Value = -5057.1667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H111").Value = 62000
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 62000
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 62000
$ws.Range("N111").Value = -70180

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1181.5
$ws.Range("I122").Value = 1067.8
$ws.Range("J122").Value = 1750
$ws.Range("K122").Value = 3203.4
$ws.Range("L122").Value = 5250
$ws.Range("M122").Value = -753.3999999999996
$ws.Range("N122").Value = -10150

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5999.75
$ws.Range("I16").Value = 3999.5
$ws.Range("J16").Value = 8000
$ws.Range("K16").Value = 3999.5
$ws.Range("L16").Value = 8000
$ws.Range("M16").Value = -3712.5
$ws.Range("N16").Value = -8574

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2933
$ws.Range("I58").Value = 2933
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 2933
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -2730

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 5999.75
$ws.Range("I113").Value = 3999.5
$ws.Range("J113").Value = 8000
$ws.Range("K113").Value = 3999.5
$ws.Range("L113").Value = 8000
$ws.Range("M113").Value = -1829.5
$ws.Range("N113").Value = -12340

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2933
$ws.Range("I136").Value = 2933
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 8799
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -6249

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2779.75
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 2779.75
$ws.Range("K68").Value = 0
$ws.Range("L68").ClearContents()
$ws.Range("M68").Value = 8339.25
$ws.Range("N68").Value = -9961.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 2779.75
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 2779.75
$ws.Range("K71").Value = 0
$ws.Range("L71").ClearContents()
$ws.Range("M71").Value = 25017.75
$ws.Range("N71").Value = -33129.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 2050
$ws.Range("I86").Value = 2050
$ws.Range("J86").Value = 2050
$ws.Range("K86").Value = 6150
$ws.Range("L86").Value = 6150
$ws.Range("M86").Value = -4964
$ws.Range("N86").Value = -8522

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 2050
$ws.Range("I89").Value = 2050
$ws.Range("J89").Value = 2050
$ws.Range("K89").Value = 18450
$ws.Range("L89").Value = 18450
$ws.Range("M89").Value = -12522
$ws.Range("N89").Value = -30306

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 5580000
$ws.Range("I11").Value = 10816667
$ws.Range("J11").Value = 343333.16
$ws.Range("K11").Value = 10816667
$ws.Range("L11").Value = 343333.16
$ws.Range("M11").Value = -10816528
$ws.Range("N11").Value = -343611.16

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 5030000
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 5030000
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 5030000
$ws.Range("N33").Value = -5030504

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 3007.6
$ws.Range("I36").Value = 1769.5
$ws.Range("J36").Value = 3833
$ws.Range("K36").Value = 1769.5
$ws.Range("L36").Value = 3833
$ws.Range("M36").Value = -1284.5
$ws.Range("N36").Value = -4803

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").ClearContents()
$ws.Range("N40").Value = 0

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 56333.332
$ws.Range("I63").Value = 50000
$ws.Range("J63").Value = 59500
$ws.Range("K63").Value = 50000
$ws.Range("L63").Value = 59500
$ws.Range("M63").Value = -49314
$ws.Range("N63").Value = -60872

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H66").Value = 56333.332
$ws.Range("I66").Value = 50000
$ws.Range("J66").Value = 59500
$ws.Range("K66").Value = 150000
$ws.Range("L66").Value = 178500
$ws.Range("M66").Value = -146568
$ws.Range("N66").Value = -185364

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4000
$ws.Range("I80").Value = 4000
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 4000
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -3002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 4000
$ws.Range("I83").Value = 4000
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 20000
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -15008

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 100000
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 100000
$ws.Range("K93").Value = 0
$ws.Range("L93").ClearContents()
$ws.Range("M93").Value = 100000
$ws.Range("N93").Value = -103744

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5498.1113
$ws.Range("I102").Value = 5185.375
$ws.Range("J102").Value = 8000
$ws.Range("K102").Value = 5185.375
$ws.Range("L102").Value = 8000
$ws.Range("M102").Value = -3563.375
$ws.Range("N102").Value = -11244

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 5000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 5000
$ws.Range("N113").Value = -9340

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 8396.6
$ws.Range("I122").Value = 7995.75
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 23987.25
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -21537.25
$ws.Range("N122").Value = -34900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1040
$ws.Range("I46").ClearContents()
$ws.Range("J46").Value = 1136.6
$ws.Range("K46").Value = 798.5
$ws.Range("L46").Value = 1136.6
$ws.Range("M46").Value = -610.5
$ws.Range("N46").Value = -1512.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H56").Value = 44485.57
$ws.Range("I56").Value = 43333
$ws.Range("J56").Value = 45350
$ws.Range("K56").Value = 43333
$ws.Range("L56").Value = 45350
$ws.Range("M56").Value = -42642
$ws.Range("N56").Value = -46732

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1897.6
$ws.Range("I61").Value = 1897.6
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1897.6
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -1695.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").ClearContents()
$ws.Range("N74").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").ClearContents()
$ws.Range("N77").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 4915.3335
$ws.Range("I100").Value = 4915.3335
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 4915.3335
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -4374.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1897.6
$ws.Range("I113").Value = 1897.6
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1897.6
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = 272.4000000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").ClearContents()
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = 0

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 28500
$ws.Range("I51").Value = 28500
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 28500
$ws.Range("L51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -27990

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("M52").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 15789
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 15789
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 15789
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -16285

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H67").Value = 15789
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 15789
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 15789
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -17505

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 90000
$ws.Range("I75").Value = 90000
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 90000
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -89064

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H78").Value = 90000
$ws.Range("I78").Value = 90000
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 270000
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -265320

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 499
$ws.Range("I81").Value = 499
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 998
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = 63

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 499
$ws.Range("I84").Value = 499
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 4990
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = 314

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 11799.8
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 11799.8
$ws.Range("K113").Value = 0
$ws.Range("L113").ClearContents()
$ws.Range("M113").Value = 35399.39999999999
$ws.Range("N113").Value = -39739.39999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2801.6
$ws.Range("I122").Value = 2502
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 7506
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -5056
$ws.Range("N122").Value = -16900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 89999
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 89999
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 89999
$ws.Range("N135").Value = -100139
